$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4935
$ws.Range("K3").Value = 8183
$ws.Range("L3").Value = 5310
$ws.Range("L4").Value = 1298
$ws.Range("L5").Value = 313
$ws.Range("L6").Value = 4482
$ws.Range("K7").Value = 27574
$ws.Range("L7").Value = 16338

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 135
$ws.Range("L7").Value = 532
$ws.Range("L8").Value = 1088
$ws.Range("L9").Value = 94
$ws.Range("L11").Value = 266
$ws.Range("L19").Value = 448
$ws.Range("L20").Value = 408
$ws.Range("L29").Value = 895
$ws.Range("L33").Value = 749
$ws.Range("L37").Value = 617
$ws.Range("L41").Value = 73
$ws.Range("L42").Value = 534
$ws.Range("L47").Value = 112
$ws.Range("L50").Value = 83
$ws.Range("L52").Value = 330
$ws.Range("L54").Value = 348
$ws.Range("L55").Value = 157
$ws.Range("K63").Value = 170
$ws.Range("L63").Value = 47
$ws.Range("L64").Value = 109
$ws.Range("L65").Value = 321
$ws.Range("L67").Value = 570
$ws.Range("L68").Value = 53
$ws.Range("L78").Value = 209
$ws.Range("L79").Value = 432
$ws.Range("L83").Value = 359
$ws.Range("L85").Value = 833
$ws.Range("L89").Value = 237
$ws.Range("L90").Value = 164
$ws.Range("L95").Value = 225
$ws.Range("L96").Value = 186
$ws.Range("L99").Value = 283
$ws.Range("K101").Value = 27574
$ws.Range("L101").Value = 16338

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 177
$ws.Range("L7").Value = 532

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 371
$ws.Range("L6").Value = 281
$ws.Range("L7").Value = 1088

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 27
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 100
$ws.Range("L6").Value = 62
$ws.Range("L7").Value = 266

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 159
$ws.Range("L7").Value = 448

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 126
$ws.Range("L7").Value = 408

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 271
$ws.Range("L6").Value = 228
$ws.Range("L7").Value = 895

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 260
$ws.Range("L7").Value = 749

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 208
$ws.Range("L7").Value = 617

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 177
$ws.Range("L4").Value = 41
$ws.Range("L6").Value = 149
$ws.Range("L7").Value = 534

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 43
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L4").Value = 23
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 330

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 63
$ws.Range("L3").Value = 84
$ws.Range("L6").Value = 170
$ws.Range("L7").Value = 348

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 157

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 29
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 117
$ws.Range("L7").Value = 321

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 168
$ws.Range("L3").Value = 218
$ws.Range("L5").Value = 14
$ws.Range("L7").Value = 570

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 67
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 142
$ws.Range("L6").Value = 96
$ws.Range("L7").Value = 432

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 84
$ws.Range("L7").Value = 359

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 247
$ws.Range("L7").Value = 833

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 71
$ws.Range("L4").Value = 38
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 237

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 89
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 56
$ws.Range("L6").Value = 53
$ws.Range("L7").Value = 186

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 77
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 283
